# "updating multisensory stimulus generation"
#
# Adds a new "cool" category to the table: a new column (mirroring the
# existing "no heat" / "heat" columns) and a new row (mirroring the
# existing "no heat" / "heat" rows), plus updates the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (pushes the "legal trials" ref/test block from
# K:L to L:M, carrying over formatting/merge automatically).
$ws.Range("J1").EntireColumn.Insert()

# Match column J's width to column I's (both width 9).
$ws.Columns.Item(9).ColumnWidth = 8.1666666666667
$ws.Columns.Item(10).ColumnWidth = 8.1666666666667

# New column header "cool" (row 2), same style as the other headers.
$ws.Range("J2").Value = "cool"

# New row 10 "cool" label + its data-entry cells, copying the formatting
# from the row above ("heat", row 9).
$ws.Range("B9:G9").Copy()
$ws.Range("B10:G10").PasteSpecial(-4122)
$ws.Range("B10").Value = "cool"

# Move the saved selection to G10.
$ws.Range("G10").Select()
